$d = $word.ActiveDocument

# Locate the paragraph containing "LOB1024: Mecânica (Requisito fraco)" -
# the three paragraphs that follow it (a blank paragraph, the
# "Ver no Jupiter..." paragraph, and the "© 2020 ..." copyright paragraph)
# are the ones that need to be removed, leaving the trailing blank
# paragraph (the one before the page-break paragraph) intact.

$lobRange = $d.Content
$lobRange.Find.Execute("LOB1024: Mecânica (Requisito fraco)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$lobRange.Expand(4)  # wdParagraph - expand to include the paragraph mark
$deleteStart = $lobRange.End

$copyrightRange = $d.Content
$copyrightRange.Find.Execute("© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$copyrightRange.Expand(4)  # wdParagraph - expand to include the paragraph mark
$deleteEnd = $copyrightRange.End

$deleteRange = $d.Range($deleteStart, $deleteEnd)
$deleteRange.Delete()
